$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update acceptedNameUsageID (col L) for the two "review" synonym rows ---
# Row 2: Acediopsylla simplex (synonym) - acceptedNameUsageID was "Lewis2404", now "Lewis12"
$ws.Range("L2").Value = "Lewis12"
# Row 3: Cediopsylla simplex (synonym) - acceptedNameUsageID was "Lewis2404", now "Lewis12"
$ws.Range("L3").Value = "Lewis12"

# --- Update TPTID (col I) values for all data rows ---
$ws.Range("I2").Value = 2790
$ws.Range("I3").Value = 3171
$ws.Range("I4").Value = 330
$ws.Range("I5").Value = 319
$ws.Range("I6").Value = 367
$ws.Range("I7").Value = 399
$ws.Range("I8").Value = 417

# --- Widen column R (scientificName, column 18) ---
$ws.Columns.Item(18).ColumnWidth = 33.5

# --- Update the active selection (no more scrolled topLeftCell, select A7) ---
$ws.Range("A7").Select()

# --- Refresh the sort state so it reflects the real data range (A2:AN8) instead of stale A2:AN62 ---
$sortRange = $ws.Range("A1:AN8")
$sortKey = $ws.Range("AN2:AN8")
$ws.Sort.SortFields.Clear()
$ws.Sort.SortFields.Add($sortKey)
$ws.Sort.SetRange($sortRange)
$ws.Sort.Header = 1
$ws.Sort.Apply()
